# Add a "2020" column (Q) to the suicide-mortality-rate table, mirroring the
# existing yearly columns (D..P), and leave the sheet with column T selected
# (whole-column selection), matching the saved view state in the target file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (borders/font/number format/alignment) of the adjacent
# "2019" column (P) onto the new "2020" column (Q) so every row keeps the
# same look (header style, body rows, bottom-bordered last row, ...).
$ws.Range("P4:P14").Copy() | Out-Null
$ws.Range("Q4").PasteSpecial(-4122) | Out-Null   # -4122 = xlPasteFormats
$excel.CutCopyMode = 0

# New "2020" year header.
$ws.Range("Q4").Value = 2020

# New data values for 2020, one per region row.
$ws.Range("Q5").Value = 4.5999999999999996
$ws.Range("Q6").Value = 4.2
$ws.Range("Q7").Value = 1.3
$ws.Range("Q8").Value = 10.8
$ws.Range("Q9").Value = 6.5
$ws.Range("Q10").Value = 2.9
$ws.Range("Q11").Value = 2.6
$ws.Range("Q12").Value = 13.1
$ws.Range("Q13").Value = 1
$ws.Range("Q14").Value = 1.3

# Match the saved selection state captured in the target workbook (column T,
# the first fully empty column after the new data, selected end-to-end).
$ws.Range("T1:T1048576").Select() | Out-Null

$wb.Save()
